# Update scripts wuth new tpm
# - Sending cluster "Neutrophils" -> "ECs"
# - Recomputed Ligand/Receptor/Edge expression metrics for row 2
# - Removed the now-redundant Neutrophils -> Neutrophils self-loop row (old row 3)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old row 3 (Neutrophils -> Neutrophils self edge), no longer present
# with the updated TPM data.
$ws.Rows.Item(3).Delete()

# Rename the sending/target cluster label used in row 2.
$ws.Range("A2").Value = "ECs"

# Refresh the expression-derived metrics for the remaining Galp -> Galr3 edge
# (Sending cluster ECs -> Target cluster FAPs) with the new TPM values.
$ws.Range("F2").Value = 0.5
$ws.Range("G2").Value = 0.1314505
$ws.Range("H2").Value = 0.262901
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.001092704189666667
$ws.Range("R2").Value = 0.006556225137999999
$ws.Range("S2").Value = 1
$ws.Range("T2").Value = 1
